$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D and E columns so numeric-looking strings are not
# auto-converted to numbers (they are stored as text in the source data).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.057.08"
$ws.Range("E2").Value = "  +1.82%  "

$ws.Range("D3").Value = "1.820.42"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").Value = "313.74"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "0.4303"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.07272"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "2.201.03"
$ws.Range("E10").Value = "  +25.42%  "

$ws.Range("D11").Value = "0.8682"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").Value = "21.26"
$ws.Range("E12").Value = "  +3.45%  "

$ws.Range("D13").Value = "5.418"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").Value = "6.626"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").Value = "0.06995"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "81.33"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "0.000008941"
$ws.Range("E18").Value = "  +2.15%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "15.29"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "27.089.68"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "2.397.49"
$ws.Range("E23").Value = "  +21.30%  "

$ws.Range("D24").Value = "11.03"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").Value = "154.33"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").Value = "1.886"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").Value = "5.245"
$ws.Range("E28").Value = "  +2.30%  "

$ws.Range("D29").Value = "1.905"
$ws.Range("E29").Value = "  +6.35%  "

$ws.Range("D30").Value = "114.71"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").Value = "0.08968"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").Value = "1.199"
$ws.Range("E32").Value = "  +6.48%  "

$ws.Range("D33").Value = "0.7560"
$ws.Range("E33").Value = "  +3.25%  "

$ws.Range("D34").Value = "4.442"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").Value = "2.804"
$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("D36").Value = "1.007"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").Value = "1.135"
$ws.Range("E37").Value = "  +4.40%  "

$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").Value = "0.01929"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("E40").Value = "  +3.06%  "

$ws.Range("D41").Value = "2.749"
$ws.Range("E41").Value = "  +5.61%  "

$ws.Range("D42").Value = "0.1654"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("D43").Value = "6.485"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("D44").Value = "8.375"
$ws.Range("E44").Value = "  +3.27%  "

$ws.Range("D45").Value = "107.18"
$ws.Range("E45").Value = "  +1.49%  "

$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "0.4590"
$ws.Range("E48").Value = "  +1.35%  "

# Rows 49-51: coin ranking reordering (RenderToken moved up, NEARProtocol and
# Cronos shift down one position each), plus updated price/volume data.
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.882"
$ws.Range("E49").Value = "  +6.08%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.649"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06226"
$ws.Range("E51").Value = "  -0.01%  "
